$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C22").Value = 99.7
$ws.Range("I22").Value = 4.9

$ws.Range("C23").Value = 99.7
$ws.Range("I23").Value = 4.5

$ws.Range("C24").Value = 100
$ws.Range("I24").Value = 4.5
